# Insert a new weekly record at row 63 ("Fruta / hortaliza, semanal"):
# this pushes the existing rows 63-95 down to 64-96 (dimension grows from
# A1:R95 to A1:R96) and the newly-opened row 63 is populated with the
# latest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63..95 down to 64..96, leaving row 63 empty and ready to fill.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the latest record.
$ws.Cells.Item(63, 1).Value = 8
$ws.Cells.Item(63, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(63, 3).Value = "Coquimbo"
$ws.Cells.Item(63, 4).Value = 44680
$ws.Cells.Item(63, 5).Value = 4
$ws.Cells.Item(63, 6).Value = 100112030
$ws.Cells.Item(63, 7).Value = "Poroto granado"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 360
$ws.Cells.Item(63, 11).Value = 27000
$ws.Cells.Item(63, 12).Value = 28000
$ws.Cells.Item(63, 13).Value = 27500
$ws.Cells.Item(63, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(63, 16).Value = 1100
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
